$d = $word.ActiveDocument

# Locate the "Credits placeholder." paragraph.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.Contains("Credits placeholder.")) {
        $targetPara = $para
        break
    }
}

# --- 1. Insert a new "Acknowledgments" Heading1 paragraph (with bookmark)
#        right before the Credits paragraph.
$targetPara.Range.InsertParagraphBefore()

$newHeadingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.Contains("Credits placeholder.")) {
        $newHeadingPara = $d.Paragraphs.Item($i - 1)
        break
    }
}

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:bookmarkStart w:id="22" w:name="acknowledgments"/><w:r><w:t xml:space="preserve">Acknowledgments</w:t></w:r><w:bookmarkEnd w:id="22"/></w:p>'
$newHeadingPara.Range.InsertXML($xml)

# --- 2. Replace the Credits placeholder text with the real credits text.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.Contains("Credits placeholder.")) {
        $para.Range.Text = "Some materials included in this export came from the following casebooks."
        break
    }
}

# --- 3. Materialize the (previously empty/default) section page setup.
$sec = $d.Sections.Item(1)
$sec.PageSetup.PageWidth = 612
$sec.PageSetup.PageHeight = 792
$sec.PageSetup.TopMargin = 72
$sec.PageSetup.BottomMargin = 72
$sec.PageSetup.LeftMargin = 72
$sec.PageSetup.RightMargin = 72
$sec.PageSetup.HeaderDistance = 36
$sec.PageSetup.FooterDistance = 36
$sec.PageSetup.Gutter = 0
$sec.PageSetup.TextColumns.Spacing = 36
